$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NIG(0.8848682078058188, 0.625848346745812, 1.2595050994395636, 3.0000016488317316)"
$ws.Range("C2").Value = "JSU(-0.9068709387578091, 0.9978708045141396, 4.489548238032727, 4.3138346885107435)"
$ws.Range("D2").Value = "JSU(-0.7340419998335979, 0.9965174489519537, 0.8526247497780218, 2.2794816092629144)"
$ws.Range("E2").Value = "NIG(0.9370527748310995, 0.7285917165562444, 4.3471121180650485, 5.442789124624513)"
